$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 16.74135869433517
$ws.Range("C2").Value2 = 10.40645936437748
$ws.Range("D2").Value2 = 9.655582376328848
$ws.Range("E2").Value2 = 13.86504579679892
$ws.Range("F2").Value2 = 30.51856598266524
$ws.Range("J2").Value2 = 9.830610554194555
$ws.Range("O2").Value2 = 22.67541293114011

$ws.Range("B3").Value2 = 16.03378893237481
$ws.Range("C3").Value2 = 9.7698729958936
$ws.Range("D3").Value2 = 9.586803596302058
$ws.Range("E3").Value2 = 13.79737846091945
$ws.Range("F3").Value2 = 30.61092083064113
$ws.Range("J3").Value2 = 9.83926703291848
$ws.Range("O3").Value2 = 22.80050043012294

$ws.Range("B4").Value2 = 15.58366149769168
$ws.Range("C4").Value2 = 9.35614875790977
$ws.Range("D4").Value2 = 9.545502157929882
$ws.Range("E4").Value2 = 13.75838719220554
$ws.Range("F4").Value2 = 30.67833523616066
$ws.Range("J4").Value2 = 9.846274117652253
$ws.Range("O4").Value2 = 22.88484312272928

$ws.Range("B5").Value2 = 15.39652788385823
$ws.Range("C5").Value2 = 9.181854375672174
$ws.Range("D5").Value2 = 9.528917775186969
$ws.Range("E5").Value2 = 13.74315272622164
$ws.Range("F5").Value2 = 30.70848496132109
$ws.Range("J5").Value2 = 9.849554982984618
$ws.Range("O5").Value2 = 22.92109894140553

$ws.Range("B6").Value2 = 15.3652384445491
$ws.Range("C6").Value2 = 9.152570486922871
$ws.Range("D6").Value2 = 9.526179189184147
$ws.Range("E6").Value2 = 13.7406629323899
$ws.Range("F6").Value2 = 30.71365257970415
$ws.Range("J6").Value2 = 9.850125459237221
$ws.Range("O6").Value2 = 22.92723277414742

$ws.Range("B7").Value2 = 15.58115240270027
$ws.Range("C7").Value2 = 9.35382115979866
$ws.Range("D7").Value2 = 9.545277481460367
$ws.Range("E7").Value2 = 13.7581790685962
$ws.Range("F7").Value2 = 30.67873102398582
$ws.Range("J7").Value2 = 9.846316642171828
$ws.Range("O7").Value2 = 22.88532446059581

$ws.Range("B8").Value2 = 16.50077614483935
$ws.Range("C8").Value2 = 10.19173382986074
$ws.Range("D8").Value2 = 9.631682356991778
$ws.Range("E8").Value2 = 13.84119083257699
$ws.Range("F8").Value2 = 30.54817883362047
$ws.Range("J8").Value2 = 9.833244147035334
$ws.Range("O8").Value2 = 22.71697242902238

$ws.Range("B9").Value2 = 18.17079458782201
$ws.Range("C9").Value2 = 11.6519965199913
$ws.Range("D9").Value2 = 9.807920635205056
$ws.Range("E9").Value2 = 14.02370279445573
$ws.Range("F9").Value2 = 30.37774723420632
$ws.Range("J9").Value2 = 9.821032360781011
$ws.Range("O9").Value2 = 22.4470980520567

$ws.Range("B10").Value2 = 19.30636051691145
$ws.Range("C10").Value2 = 12.61204613125857
$ws.Range("D10").Value2 = 9.940761717899367
$ws.Range("E10").Value2 = 14.16901682208172
$ws.Range("F10").Value2 = 30.30547356949262
$ws.Range("J10").Value2 = 9.820236166297816
$ws.Range("O10").Value2 = 22.28616836538848

$ws.Range("B11").Value2 = 19.80146181603468
$ws.Range("C11").Value2 = 13.02417847875576
$ws.Range("D11").Value2 = 10.0017521746106
$ws.Range("E11").Value2 = 14.23737480410113
$ws.Range("F11").Value2 = 30.28422967106811
$ws.Range("J11").Value2 = 9.821645231203814
$ws.Range("O11").Value2 = 22.22119525824397

$ws.Range("B12").Value2 = 19.98574531011602
$ws.Range("C12").Value2 = 13.17670081036791
$ws.Range("D12").Value2 = 10.02491357235228
$ws.Range("E12").Value2 = 14.26356765386705
$ws.Range("F12").Value2 = 30.2778673695386
$ws.Range("J12").Value2 = 9.8224329212455
$ws.Range("O12").Value2 = 22.19778618970858

$ws.Range("B13").Value2 = 19.94620052110497
$ws.Range("C13").Value2 = 13.14401003402625
$ws.Range("D13").Value2 = 10.01992267020315
$ws.Range("E13").Value2 = 14.25791316902726
$ws.Range("F13").Value2 = 30.27916265004493
$ws.Range("J13").Value2 = 9.822251987524458
$ws.Range("O13").Value2 = 22.20277444202607

$ws.Range("B14").Value2 = 19.81668765893334
$ws.Range("C14").Value2 = 13.03679761370037
$ws.Range("D14").Value2 = 10.00365644707172
$ws.Range("E14").Value2 = 14.23952364845335
$ws.Range("F14").Value2 = 30.28367247621379
$ws.Range("J14").Value2 = 9.821704946136114
$ws.Range("O14").Value2 = 22.2192453730326

$ws.Range("B15").Value2 = 19.73693742919215
$ws.Range("C15").Value2 = 12.970665371858
$ws.Range("D15").Value2 = 9.993701019037966
$ws.Range("E15").Value2 = 14.22829902636569
$ws.Range("F15").Value2 = 30.28665421433966
$ws.Range("J15").Value2 = 9.82140293981368
$ws.Range("O15").Value2 = 22.22949021188955

$ws.Range("B16").Value2 = 19.27356233458898
$ws.Range("C16").Value2 = 12.58461671118618
$ws.Range("D16").Value2 = 9.936785847967865
$ws.Range("E16").Value2 = 14.16459338536217
$ws.Range("F16").Value2 = 30.30709684357664
$ws.Range("J16").Value2 = 9.820179674855895
$ws.Range("O16").Value2 = 22.29058113146246

$ws.Range("B17").Value2 = 18.9837108840617
$ws.Range("C17").Value2 = 12.34148232660939
$ws.Range("D17").Value2 = 9.902003145956781
$ws.Range("E17").Value2 = 14.12607743442113
$ws.Range("F17").Value2 = 30.32262454306343
$ws.Range("J17").Value2 = 9.819882520788809
$ws.Range("O17").Value2 = 22.33017518539728

$ws.Range("B18").Value2 = 18.81498034558329
$ws.Range("C18").Value2 = 12.19932252460664
$ws.Range("D18").Value2 = 9.882050701333373
$ws.Range("E18").Value2 = 14.10413727660702
$ws.Range("F18").Value2 = 30.3326502161649
$ws.Range("J18").Value2 = 9.819878368973054
$ws.Range("O18").Value2 = 22.3537230878124

$ws.Range("B19").Value2 = 18.75750862612264
$ws.Range("C19").Value2 = 12.15079225009724
$ws.Range("D19").Value2 = 9.875304814568995
$ws.Range("E19").Value2 = 14.09674584057943
$ws.Range("F19").Value2 = 30.33623241797899
$ws.Range("J19").Value2 = 9.819905621004212
$ws.Range("O19").Value2 = 22.36182868655841

$ws.Range("B20").Value2 = 19.01477555135674
$ws.Range("C20").Value2 = 12.36760408186228
$ws.Range("D20").Value2 = 9.905700379962653
$ws.Range("E20").Value2 = 14.13015557896135
$ws.Range("F20").Value2 = 30.32085823847469
$ws.Range("J20").Value2 = 9.819896897863842
$ws.Range("O20").Value2 = 22.32588009995824

$ws.Range("B21").Value2 = 19.85481640185666
$ws.Range("C21").Value2 = 13.06838468946341
$ws.Range("D21").Value2 = 10.00843256910418
$ws.Range("E21").Value2 = 14.24491689763332
$ws.Range("F21").Value2 = 30.28230210609031
$ws.Range("J21").Value2 = 9.821858734463939
$ws.Range("O21").Value2 = 22.21437494928539

$ws.Range("B22").Value2 = 20.38513559736465
$ws.Range("C22").Value2 = 13.505732636694
$ws.Range("D22").Value2 = 10.07594961684785
$ws.Range("E22").Value2 = 14.32170224031981
$ws.Range("F22").Value2 = 30.26691251689434
$ws.Range("J22").Value2 = 9.824621670701292
$ws.Range("O22").Value2 = 22.14846915645739

$ws.Range("B23").Value2 = 20.10383777703414
$ws.Range("C23").Value2 = 13.27420243844974
$ws.Range("D23").Value2 = 10.03988507414146
$ws.Range("E23").Value2 = 14.28056311540914
$ws.Range("F23").Value2 = 30.27422597261683
$ws.Range("J23").Value2 = 9.823011781666294
$ws.Range("O23").Value2 = 22.1830031364816

$ws.Range("B24").Value2 = 19.00073771978766
$ws.Range("C24").Value2 = 12.35580184474569
$ws.Range("D24").Value2 = 9.90402872044559
$ws.Range("E24").Value2 = 14.1283112156842
$ws.Range("F24").Value2 = 30.32165336338373
$ws.Range("J24").Value2 = 9.819889878731043
$ws.Range("O24").Value2 = 22.32781946435464

$ws.Range("B25").Value2 = 17.73442768632358
$ws.Range("C25").Value2 = 11.27670342595027
$ws.Range("D25").Value2 = 9.759595803215673
$ws.Range("E25").Value2 = 13.97229656838425
$ws.Range("F25").Value2 = 30.41460723276328
$ws.Range("J25").Value2 = 9.822899182246774
$ws.Range("O25").Value2 = 22.51359136460945

